# Update latest output (run 184)
# Applies cell value changes to the Schedule and Detailed sheets

$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# --- Schedule sheet updates ---
$wsSchedule.Range("B2").Value = 46071.27083333334
$wsSchedule.Range("C2").Value = 6.5
$wsSchedule.Range("D2").Value = 24.57
$wsSchedule.Range("E2").Value = 926.9015145
$wsSchedule.Range("F2").Value = 37.72492936507936
$wsSchedule.Range("A4").Value = 46071.95833333334
$wsSchedule.Range("C4").Value = 4.5
$wsSchedule.Range("D4").Value = 17.01
$wsSchedule.Range("E4").Value = 710.3674597499999
$wsSchedule.Range("F4").Value = 41.76175542328042
$wsSchedule.Range("E5").Value = 715.6762470000002
$wsSchedule.Range("F5").Value = 22.27439299719889

# --- Detailed sheet updates ---
$wsDetailed.Range("E14").Value = "ON"
$wsDetailed.Range("B38").Value = 56.57435
$wsDetailed.Range("B39").Value = 67.5433
$wsDetailed.Range("B40").Value = 157.27893
$wsDetailed.Range("C40").Value = "historical"
$wsDetailed.Range("B41").Value = 154.2
$wsDetailed.Range("C41").Value = "historical"
$wsDetailed.Range("B42").Value = 146.17011
$wsDetailed.Range("C42").Value = "historical"
$wsDetailed.Range("B43").Value = 143.63034
$wsDetailed.Range("C43").Value = "historical"
$wsDetailed.Range("B44").Value = 138.42
$wsDetailed.Range("C44").Value = "historical"
$wsDetailed.Range("B45").Value = 108.89
$wsDetailed.Range("C45").Value = "historical"
$wsDetailed.Range("B46").Value = 108.89
$wsDetailed.Range("C46").Value = "historical"
$wsDetailed.Range("B47").Value = 95.11075
$wsDetailed.Range("C47").Value = "historical"
$wsDetailed.Range("E47").Value = "OFF"
$wsDetailed.Range("B48").Value = 86.06352
$wsDetailed.Range("C48").Value = "historical"
$wsDetailed.Range("B49").Value = 73.20005999999999
$wsDetailed.Range("C49").Value = "historical"
$wsDetailed.Range("B50").Value = 80.02367
$wsDetailed.Range("B51").Value = 84.70398
$wsDetailed.Range("B52").Value = 84.79000000000001
$wsDetailed.Range("B53").Value = 79.95038
$wsDetailed.Range("B54").Value = 79.9504
$wsDetailed.Range("B58").Value = 79.95
$wsDetailed.Range("B59").Value = 83.47329999999999
$wsDetailed.Range("B60").Value = 90.54376000000001
$wsDetailed.Range("B61").Value = 98.17335
$wsDetailed.Range("B62").Value = 106.84888
$wsDetailed.Range("B63").Value = 108.89
$wsDetailed.Range("B64").Value = 100.01
$wsDetailed.Range("B65").Value = 60.51663
$wsDetailed.Range("B66").Value = 57.06
$wsDetailed.Range("B67").Value = 51.29678
$wsDetailed.Range("B68").Value = 36.06
$wsDetailed.Range("B70").Value = 36.06
$wsDetailed.Range("B74").Value = 36.06
$wsDetailed.Range("B77").Value = 54.35344
$wsDetailed.Range("B78").Value = 36.07
$wsDetailed.Range("B80").Value = 57.06
$wsDetailed.Range("B82").Value = 64.9019
$wsDetailed.Range("B83").Value = 60.0036
$wsDetailed.Range("B84").Value = 64.12925
$wsDetailed.Range("B85").Value = 64.22445999999999
$wsDetailed.Range("B86").Value = 63.53876
$wsDetailed.Range("B87").Value = 135.50041
$wsDetailed.Range("B88").Value = 166.61417
$wsDetailed.Range("B89").Value = 156.49832
$wsDetailed.Range("B90").Value = 180.12572
$wsDetailed.Range("B91").Value = 138.60419
$wsDetailed.Range("B93").Value = 104.74819
$wsDetailed.Range("B94").Value = 92.70107
$wsDetailed.Range("B95").Value = 92.65982
$wsDetailed.Range("B97").Value = 85.65000000000001
